$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.63 = 6024.39 pesos`n✅ 6024.39 pesos = 1.62 = 948.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 615
$ws2.Range("O10").Value = 3705
$ws2.Range("N12").Value = 3729.99
$ws2.Range("O12").Value = 587.05
